$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.196.52"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.24%  "

$ws.Range("D3").Value = "'2.057.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.81%  "

$ws.Range("D5").Value = "'231.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.25%  "

$ws.Range("D6").Value = "'0.616"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.93%  "

$ws.Range("D8").Value = "'57.54"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.17%  "

$ws.Range("D9").Value = "'0.381"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.11%  "

$ws.Range("D10").Value = "'57.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.75%  "

$ws.Range("D11").Value = "'0.0756"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.53%  "

$ws.Range("E12").Value = "  +1.26%  "

$ws.Range("D13").Value = "'2.364.74"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.70%  "

$ws.Range("D14").Value = "'14.30"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.73%  "

$ws.Range("D15").Value = "'20.67"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.37%  "

$ws.Range("D16").Value = "'0.771"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.90%  "

$ws.Range("E17").Value = "  +1.25%  "

$ws.Range("D18").Value = "'2.061.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.67%  "

$ws.Range("D19").Value = "'37.204.04"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.64%  "

$ws.Range("D20").Value = "'6.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +15.81%  "

$ws.Range("D21").Value = "'68.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.14%  "

$ws.Range("D22").Value = "'0.0₃0806"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.83%  "

$ws.Range("D23").Value = "'224.40"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.73%  "

$ws.Range("E24").Value = "  +0.04%  "

$ws.Range("E25").Value = "  +2.51%  "

$ws.Range("E26").Value = "  +0.81%  "

$ws.Range("D27").Value = "'165.41"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.90%  "

$ws.Range("E28").Value = "  +8.09%  "

$ws.Range("D29").Value = "'8.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.94%  "

$ws.Range("E30").Value = "  -0.78%  "

$ws.Range("D31").Value = "'18.98"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.80%  "

$ws.Range("E32").Value = "  -0.05%  "

$ws.Range("D33").Value = "'4.45"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.20%  "

$ws.Range("D34").Value = "'0.0611"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.24%  "

$ws.Range("E35").Value = "  +3.48%  "

$ws.Range("E36").Value = "  +6.83%  "

$ws.Range("E37").Value = "  +0.04%  "

$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").Value = "'1.74"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.08%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'3.26"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.11%  "

$ws.Range("B40").Value = "THORChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D40").Value = "'5.79"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.30%  "

$ws.Range("D41").Value = "'4.63"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +17.11%  "

$ws.Range("D42").Value = "'2.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.74%  "

$ws.Range("D43").Value = "'1.489.26"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.15%  "

$ws.Range("D44").Value = "'96.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.10%  "

$ws.Range("E45").Value = "  +5.12%  "

$ws.Range("D46").Value = "'0.0925"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.06%  "

$ws.Range("E47").Value = "  +3.40%  "

$ws.Range("D48").Value = "'15.28"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.06%  "

$ws.Range("E49").Value = "  +2.35%  "

$ws.Range("B50").Value = "MXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D50").Value = "'2.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.55%  "

$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").Value = "'7.14"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.36%  "
